# Add a new "Contributions of Team Members" section at the end of the
# document body, after the paragraph that currently ends the document
# ("...Algorithm could address this concern, particularly one which
# encouraged speciation.").

$d = $word.ActiveDocument

# The paragraph we need to insert after is the very last paragraph in
# the main story (just before the final section break).
$lastPara = $d.Paragraphs.Last

# --- New Heading 1 paragraph: "Contributions of Team Members" ---
$lastPara.Range.InsertParagraphAfter()
$headingPara = $d.Paragraphs($d.Paragraphs.Count)
$headingPara.Range.Text = "Contributions of Team Members"
$headingPara.Style = "Heading 1"

# --- New plain paragraph: Randall's contributions ---
$headingPara.Range.InsertParagraphAfter()
$randallPara = $d.Paragraphs($d.Paragraphs.Count)
$randallPara.Style = "Normal"
$randallPara.Range.Text = "Randall " + [char]0x2013 + " algorithm development, algorithm implementation"

# --- New plain paragraph: Conor's contributions ---
$randallPara.Range.InsertParagraphAfter()
$conorPara = $d.Paragraphs($d.Paragraphs.Count)
$conorPara.Style = "Normal"
$conorPara.Range.Text = "Conor " + [char]0x2013 + " algorithm development, experimentation"
